# Simulated Wild Card round and logged it
#
# Updates the Steelers "Players Data" workbook with the results of the
# simulated Wild Card playoff game:
#  - "Rushing" sheet: carry/yard totals bumped for the backs who ran the ball
#  - "Receiving" sheet: target/catch totals bumped for everyone who played,
#    and a new receiver (J.Smith-Schuster) who caught his first pass of the
#    game is inserted into the table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Rushing
# ---------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# N.Harris (row 4)
$rushing.Range("C4").Value = 185
$rushing.Range("D4").Value = 110

# B.Snell (row 5)
$rushing.Range("C5").Value = 22
$rushing.Range("D5").Value = 13
$rushing.Range("E5").Value = 4

# K.Ballage (row 6)
$rushing.Range("C6").Value = 10
$rushing.Range("D6").Value = 4
$rushing.Range("E6").Value = 2

# C.Claypool (row 10)
$rushing.Range("C10").Value = 7
$rushing.Range("D10").Value = 6

# ---------------------------------------------------------------------
# Sheet 2: Receiving
# ---------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# N.Harris (row 2)
$receiving.Range("C2").Value = 93
$receiving.Range("D2").Value = 73

# B.Snell (row 3)
$receiving.Range("C3").Value = 6
$receiving.Range("D3").Value = 4

# D.Watt (row 6)
$receiving.Range("C6").Value = 3

# D.Johnson (row 7)
$receiving.Range("C7").Value = 142
$receiving.Range("D7").Value = 99
$receiving.Range("E7").Value = 42
$receiving.Range("G7").Value = 22
$receiving.Range("H7").Value = 12

# J.Smith-Schuster caught his first pass of the game - insert him as a new
# row right after D.Johnson (row 7), pushing everyone below down by one.
$receiving.Rows(8).Insert()

# Copy row formatting down from the row above so the new row matches the
# rest of the table (bold/bordered index cell, plain data cells).
$receiving.Range("A7:H7").Copy()
$receiving.Range("A8:H8").PasteSpecial(-4122)

$receiving.Range("B8").Value = "J.Smith-Schuster"
$receiving.Range("C8").Value = 7
$receiving.Range("D8").Value = 5
$receiving.Range("E8").Value = 1
$receiving.Range("F8").Value = 0
$receiving.Range("G8").Value = 0
$receiving.Range("H8").Value = 0

# C.Claypool, now on row 9
$receiving.Range("C9").Value = 81
$receiving.Range("D9").Value = 50
$receiving.Range("E9").Value = 33
$receiving.Range("F9").Value = 13
$receiving.Range("G9").Value = 12

# J.Washington, now on row 10
$receiving.Range("C10").Value = 34
$receiving.Range("D10").Value = 22
$receiving.Range("E10").Value = 14
$receiving.Range("F10").Value = 5
$receiving.Range("G10").Value = 9
$receiving.Range("H10").Value = 6

# R.McCloud, now on row 11
$receiving.Range("C11").Value = 58
$receiving.Range("D11").Value = 38
$receiving.Range("E11").Value = 10
$receiving.Range("F11").Value = 2
$receiving.Range("G11").Value = 10
$receiving.Range("H11").Value = 3

# A.Miller, now on row 12 (unchanged stats)
$receiving.Range("C12").Value = 1
$receiving.Range("D12").Value = 1

# C.White, now on row 13 (unchanged stats)
$receiving.Range("C13").Value = 3
$receiving.Range("D13").Value = 3

# E.Ebron, now on row 14 (unchanged stats)
$receiving.Range("C14").Value = 17
$receiving.Range("D14").Value = 12
$receiving.Range("G14").Value = 3
$receiving.Range("H14").Value = 1

# P.Freiermuth, now on row 15
$receiving.Range("C15").Value = 78
$receiving.Range("D15").Value = 61

# Z.Gentry, now on row 16
$receiving.Range("C16").Value = 25
$receiving.Range("D16").Value = 21
$receiving.Range("E16").Value = 3
$receiving.Range("F16").Value = 2
$receiving.Range("G16").Value = 4
$receiving.Range("H16").Value = 2

# K.Rader, now on row 17 (unchanged stats - row just shifted down)
$receiving.Range("C17").Value = 3
$receiving.Range("D17").Value = 2

# Renumber the leading index column (0-based row counter) now that a row
# was inserted in the middle of the table.
for ($i = 2; $i -le 17; $i++) {
    $receiving.Cells.Item($i, 1).Value = $i - 2
}
